$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- sheet2 ("roll") new sweep data: rows 2-12, columns A,B,D,E,F ---
$data = @(
    @(0, 2.5, 1574.9, 214.65, 431.01),
    @(0, 2,   1574.9, 217.89, 431.31),
    @(0, 1.5, 1574.9, 220.97, 431.55),
    @(0, 1,   1575,   223.91, 431.74),
    @(0, 0.5, 1575,   226.7,  431.89),
    @(0, 0,   1575,   229.36, 432),
    @(0, -0.5,1575,   231.89, 432.07),
    @(0, -1,  1575,   234.3,  432.11),
    @(0, -1.5,1575.1, 236.57, 432.11),
    @(0, -2,  1575.1, 238.72, 432.07),
    @(0, -2.5,1575.1, 240.73, 432)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = 2 + $i
    $row = $data[$i]
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 4).Value = $row[2]
    $ws2.Cells.Item($r, 5).Value = $row[3]
    $ws2.Cells.Item($r, 6).Value = $row[4]
}

# --- delete the two now-obsolete trailing rows (13 & 14) ---
$ws2.Rows("13:14").Delete()

# --- rewrite I/J formulas for the remaining 11 data rows (2-12) ---
for ($r = 2; $r -le 12; $r++) {
    $ws2.Cells.Item($r, 9).Formula = "=F$r-`$F`$7"
}
for ($r = 2; $r -le 12; $r++) {
    $next = $r + 1
    $ws2.Cells.Item($r, 10).Formula = "=(I$r-I$next)/(B$r-B$next)"
}

# --- re-apply the (cosmetic) descending sort on column B that Excel recorded ---
$ws2.Sort.SortFields.Clear()
$ws2.Sort.SortFields.Add($ws2.Range("B1"), 0, 2, 0, 0)
$ws2.Sort.SetRange($ws2.Range("A1:F14"))
$ws2.Sort.Header = 0
$ws2.Sort.Apply()

# --- view/selection state ---
$ws1.Range("J18").Select()
$ws2.Activate()
$ws2.Range("I2:J12").Select()
